$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E146").Formula = "=0.6*C2+0.4*D2"

$ws.Range("F144").Select()
$excel.ActiveWindow.ScrollRow = 137
